# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Reorder two pairs of country names (shared-string order swap in the diff) ---
# Row 59/60: Barein <-> Kazajistan
$ws.Range("A59").Value = "Kazajistan"
$ws.Range("A60").Value = "Barein"

# Row 192/193: Belice <-> Nueva Caledonia
$ws.Range("A192").Value = "Nueva Caledonia"
$ws.Range("A193").Value = "Belice"

# --- Update numeric data values ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1368036
$ws.Range("C4").Value = 398
$ws.Range("E4").Value = 1030911
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 80789

# Row 10 - Alemania
$ws.Range("B10").Value = 171999
$ws.Range("C10").Value = 120
$ws.Range("E10").Value = 18830

# Row 11 - Brasil
$ws.Range("B11").Value = 163510
$ws.Range("C11").Value = 811
$ws.Range("E11").Value = 87346
$ws.Range("G11").Value = 84
$ws.Range("H11").Value = 11207

# Row 20 - Arabia Saudita
$ws.Range("B20").Value = 41014
$ws.Range("C20").Value = 1966
$ws.Range("D20").Value = 12737
$ws.Range("E20").Value = 28022
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 255

# Row 59 - now Kazajistan (values per diff for row 59)
$ws.Range("B59").Value = 5160
$ws.Range("C59").Value = 70
$ws.Range("D59").Value = 2020
$ws.Range("E59").Value = 3108
$ws.Range("F59").Value = 33
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 32

# Row 60 - now Barein (values per diff for row 60)
$ws.Range("B60").Value = 5157
$ws.Range("C60").Value = 216
$ws.Range("D60").Value = 2152
$ws.Range("E60").Value = 2997
$ws.Range("F60").Value = 2
$ws.Range("H60").Value = 8

# Row 75 - Uzbekistan
$ws.Range("B75").Value = 2482
$ws.Range("C75").Value = 64
$ws.Range("E75").Value = 490

# Row 118 - Georgia
$ws.Range("E118").Value = 310
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 11

# Row 192 - now Nueva Caledonia (values per diff for row 192)
$ws.Range("D192").Value = 18
$ws.Range("H192").Value = 0

# Row 193 - now Belice (values per diff for row 193)
$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2
